$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 173
$ws1.Range("F5").Value = 1312
$ws1.Range("F6").Value = 384
$ws1.Range("F8").Value = 923
$ws1.Range("F9").Value = 734
$ws1.Range("F10").Value = 206
$ws1.Range("F12").Value = 154
$ws1.Range("F15").Value = 3040
$ws1.Range("F16").Value = 2671
$ws1.Range("F21").Value = 249
$ws1.Range("F23").Value = 5444
$ws1.Range("F25").Value = 1000
$ws1.Range("F26").Value = 35
$ws1.Range("F27").Value = 65
$ws1.Range("F28").Value = 370
$ws1.Range("F29").Value = 1150
$ws1.Range("F31").Value = 80

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 41
$ws2.Range("F13").Value = 623
$ws2.Range("F23").Value = 331
$ws2.Range("F25").Value = 4005
$ws2.Range("F32").Value = 179

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2520
$ws3.Range("F6").Value = 1082
$ws3.Range("F9").Value = 1380
$ws3.Range("F10").Value = 384

# B12/E12 hold dotted text that Excel's COM layer would otherwise
# auto-convert to a date serial when assigned via .Value. Enter the text
# via a formula (so it is not re-parsed as a date) and then convert the
# formula to a static value with Copy / PasteSpecial (values only), which
# preserves the literal text and the cell's original (default) style.
$b12 = $ws3.Range("B12")
$b12.Formula = '="2024.02.24"'
$b12.Copy()
$b12.PasteSpecial(-4163)

$e12 = $ws3.Range("E12")
$e12.Formula = '="2024.02.24 19:00-03.03 20:10"'
$e12.Copy()
$e12.PasteSpecial(-4163)

$excel.CutCopyMode = 0

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2520
$ws4.Range("F9").Value = 1082
$ws4.Range("F10").Value = 1380
$ws4.Range("F11").Value = 384
$ws4.Range("F13").Value = 173
$ws4.Range("F14").Value = 1312
$ws4.Range("F15").Value = 384
$ws4.Range("F16").Value = 923
$ws4.Range("F17").Value = 734
$ws4.Range("F19").Value = 206
$ws4.Range("F21").Value = 154
$ws4.Range("F23").Value = 3040
$ws4.Range("F24").Value = 2671
$ws4.Range("F28").Value = 41
$ws4.Range("F29").Value = 249
$ws4.Range("F31").Value = 5444
$ws4.Range("F33").Value = 1000
$ws4.Range("F34").Value = 623
$ws4.Range("F35").Value = 35
$ws4.Range("F36").Value = 65
$ws4.Range("F37").Value = 370
$ws4.Range("F42").Value = 331
$ws4.Range("F47").Value = 179
$ws4.Range("F49").Value = 80
